$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from an existing header cell (AC1) onto the new header
# cells so they reuse the same style (bold, bordered header look) instead
# of creating brand new style entries.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 73
    $ws.Cells.Item($r, 32).Value = 0
}
